$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref='D2'; Value='320.55'},
    @{Ref='E2'; Value='5.80%'},
    @{Ref='D3'; Value='49.53'},
    @{Ref='E3'; Value='14.60%'},
    @{Ref='D4'; Value='5.261'},
    @{Ref='E4'; Value='3.86%'},
    @{Ref='D5'; Value='0.07993'},
    @{Ref='E5'; Value='4.03%'},
    @{Ref='D6'; Value='4.580'},
    @{Ref='E6'; Value='3.79%'},
    @{Ref='D7'; Value='1.412'},
    @{Ref='E7'; Value='35.11%'},
    @{Ref='D8'; Value='1.646'},
    @{Ref='E8'; Value='1.39%'},
    @{Ref='D9'; Value='0.1305'},
    @{Ref='E9'; Value='3.50%'},
    @{Ref='D10'; Value='0.1960'},
    @{Ref='E10'; Value='5.62%'},
    @{Ref='D11'; Value='0.09460'},
    @{Ref='E11'; Value='3.21%'},
    @{Ref='D12'; Value='0.04590'},
    @{Ref='E12'; Value='10.13%'},
    @{Ref='E13'; Value='-0.25%'},
    @{Ref='D14'; Value='0.001322'},
    @{Ref='E14'; Value='3.20%'},
    @{Ref='D15'; Value='0.04165'},
    @{Ref='E15'; Value='-0.20%'},
    @{Ref='D16'; Value='0.005905'},
    @{Ref='E16'; Value='2.79%'},
    @{Ref='E17'; Value='-0.04%'},
    @{Ref='E18'; Value='3.34%'},
    @{Ref='D19'; Value='0.3462'},
    @{Ref='E19'; Value='3.20%'},
    @{Ref='D20'; Value='8.204'},
    @{Ref='E20'; Value='-5.17%'},
    @{Ref='D21'; Value='0.1389'},
    @{Ref='E21'; Value='1.60%'},
    @{Ref='D22'; Value='0.3089'},
    @{Ref='E22'; Value='-3.36%'},
    @{Ref='D23'; Value='0.001312'},
    @{Ref='E23'; Value='2.27%'},
    @{Ref='D24'; Value='0.004263'},
    @{Ref='E24'; Value='-4.50%'},
    @{Ref='D25'; Value='0.0001347'},
    @{Ref='E25'; Value='-0.08%'},
    @{Ref='D26'; Value='0.0003530'},
    @{Ref='E26'; Value='-95.25%'},
    @{Ref='D38'; Value='0.02690'},
    @{Ref='E38'; Value='9.49%'},
    @{Ref='D39'; Value='0.06033'},
    @{Ref='E39'; Value='14.33%'},
    @{Ref='D40'; Value='0.01092'},
    @{Ref='E40'; Value='83.78%'},
    @{Ref='D41'; Value='0.008002'},
    @{Ref='E41'; Value='4.33%'},
    @{Ref='D42'; Value='0.1443'},
    @{Ref='E42'; Value='7.04%'},
    @{Ref='D43'; Value='0.007766'},
    @{Ref='E43'; Value='5.32%'},
    @{Ref='D44'; Value='0.008674'},
    @{Ref='E44'; Value='14.68%'},
    @{Ref='D45'; Value='0.3198'},
    @{Ref='E45'; Value='6.22%'},
    @{Ref='D46'; Value='0.00006610'},
    @{Ref='E46'; Value='-1.57%'},
    @{Ref='D47'; Value='0.00000000748'},
    @{Ref='E47'; Value='-0.08%'},
    @{Ref='E48'; Value='67.00%'},
    @{Ref='D49'; Value='0.003990'},
    @{Ref='E49'; Value='-5.03%'},
    @{Ref='D50'; Value='0.00002095'},
    @{Ref='E50'; Value='-0.08%'},
    @{Ref='D51'; Value='0.0001995'},
    @{Ref='E51'; Value='-0.08%'}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

Write-Output ("Applied " + $updates.Count + " updates")
